$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newCasesQuery = "MATCH (c:case)`n MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)`n MATCH (f:file)-[*]->(c)`n WHERE c.gender='MALE'`nRETURN DISTINCT`n    c.case_id AS ``Case ID``,`n     ct.clinical_trial_designation AS ``Trial Code``,`n     a.arm_id AS Arm,`n      a.arm_drug AS ``Arm Treatment``,`nc.disease AS Diagnosis,`n  c.gender AS Gender,`n    c.race AS Race,`n    c.ethnicity AS Ethnicity"

$newFilesStatQuery = "`nMATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nWITH f,a,ct,c`n    WHERE c.gender = ""MALE""`nRETURN`n    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,`n    COUNT(DISTINCT c.case_id) AS Cases,`n    COUNT(DISTINCT f) AS Files"

$ws.Range("B2").Value = $newCasesQuery
$ws.Range("C2").Value = $newFilesStatQuery

$ws.Rows.Item(2).RowHeight = 195

$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 2
